# Remove "for demo" sample entries from the "Excluded structures" sheet
# to create a blank slate for the tracking table template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Excluded structures")
$ws.Activate()

# Clear the demo row (row 2) contents while keeping formatting.
$ws.Range("A2:K2").ClearContents()

# Reset the view so that the top-left visible cell is back at the
# beginning of the sheet and the active selection is A2.
$ws.Range("A1").Select()
$ws.Range("A2").Select()
